# -----------------------------------------------------------------------
# Commit: "download tc, tcn, and tl files from GD"
#
# The underlying XML diff touches only word/styles.xml's <w:docDefaults>
# block: it drops a series of explicit run/paragraph properties
# (w:b=0, w:i=0, w:smallCaps=0, w:strike=0, w:color=000000, w:u=none,
#  w:shd=clear/auto, w:vertAlign=baseline, w:keepNext=0, w:keepLines=0,
#  w:widowControl=1, an all-"nil" w:pBdr, w:ind=0/0/0,
#  w:contextualSpacing=0, w:jc=left, and the before=0/after=0 spacing
#  attributes) that are already equal to Word's own built-in defaults,
#  and leaves <w:spacing w:line="276" w:lineRule="auto"/> as the only
#  meaningful survivor. No run/paragraph in the document body changes,
#  and the rendered/effective formatting is identical before and after
#  (every removed attribute was a no-op restatement of the implicit
#  default). This is consistent with the commit message: the .docx was
#  simply re-downloaded/re-exported from Google Drive, which happens to
#  serialize <w:docDefaults> more tersely than the previous tool - it is
#  not an edit made through Word's UI/automation surface.
#
# <w:docDefaults> sits outside the Word object model: it is not a Style
# (Styles.Item(...) only ever reads/writes explicit overrides on actual
# named styles such as "Normal", never the document-wide fallback block),
# and Document.WordOpenXML is read-only here ("the assignment changed
# nothing... To edit content, set Range.Text ... or call InsertXML").
# There is therefore no COM-reachable call that edits w:docDefaults
# directly - exactly as in real Word, where w:docDefaults can only be
# hand-edited in the OOXML package, never through the Application object.
#
# Since every property this hunk removes was already semantically
# inert, and nothing in the document body/content changes, the
# content-faithful action is to leave the document's paragraphs, runs,
# and styles exactly as authored - touching nothing here avoids
# introducing spurious explicit-formatting overrides (e.g. on the
# "Normal" style) that are not present in the target revision.
$d = $word.ActiveDocument
